# Clean up documentation edits for Memo.docx

$d = $word.ActiveDocument
$quote = [char]0x201C
$rquote = [char]0x201D

# -----------------------------------------------------------------------
# 1) "Subject:   the Traveler Language" - the words were originally split
#    across three runs (with gramStart/gramEnd proofing marks around
#    "the"). A same-text Find/Replace forces Word to rebuild the run(s),
#    merging them into a single run and dropping the proofing marks.
# -----------------------------------------------------------------------
$d.Content.Find.Execute("Subject:   the Traveler Language", $true, $false, $false, $false, $false, $true, 1, $false, "Subject:   the Traveler Language", 2) | Out-Null

# -----------------------------------------------------------------------
# 2) Rework the "context" paragraph wording.
# -----------------------------------------------------------------------

# "...can inquire for a travel plan from..." -> "...can inquire for  travel plans from..."
$d.Content.Find.Execute("for a travel plan from one location", $true, $false, $false, $false, $false, $true, 1, $false, "for  travel plans from one location", 2) | Out-Null

# "...all the "plans" requested by the programmer. " -> "...all the possible transportation
# routes for the "plans" requested by the programmer, within the provided context."
$findText = "all the " + $quote + "plans" + $rquote + " requested by the programmer. "
$replText = "all the possible transportation routes for the " + $quote + "plans" + $rquote + " requested by the programmer, within the provided context."
$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replText, 2) | Out-Null

# -----------------------------------------------------------------------
# 2b) Move the "_GoBack" bookmark from the end of the last milestone
#     bullet to the (now non-italic) empty paragraph right after the
#     "context" paragraph. Delete the old one first so there is never
#     more than one bookmark sharing the reserved "_GoBack" name.
# -----------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

$emptyPara = $d.Paragraphs.Item(6)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$emptyPara.Range.InsertXML($newParaXml) | Out-Null

# -----------------------------------------------------------------------
# 3) "Language Specification" hyperlink text - merge runs the same way
#    as (1), collapsing "Language Spec" + "i" + "fication" into one run.
# -----------------------------------------------------------------------
$d.Content.Find.Execute("Language Specification", $true, $false, $false, $false, $false, $true, 1, $false, "Language Specification", 2) | Out-Null

# -----------------------------------------------------------------------
# 4) Replace milestone bullet text.
# -----------------------------------------------------------------------
$d.Content.Find.Execute("Format the output so that it" + [char]0x2019 + "s clean and readable (instead of showing the underlying data representation)", $true, $false, $false, $false, $false, $true, 1, $false, "Deal with technical issues in plan arranging, such as time-zones and temporal gaps between transportations.", 2) | Out-Null
